$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6
$ws.Range("H6").Value = 794.375
$ws.Range("I6").Value = 225.83333
$ws.Range("K6").Value = 677.49999
$ws.Range("M6").Value = -565.49999

# Row 31
$ws.Range("H31").Value = 2337.1428
$ws.Range("I31").Value = 32
$ws.Range("K31").Value = 96
$ws.Range("M31").Value = 134

# Row 69
$ws.Range("H69").Value = 8816
$ws.Range("I69").Value = 13000
$ws.Range("K69").Value = 39000
$ws.Range("M69").Value = -38126

# Row 72
$ws.Range("H72").Value = 8816
$ws.Range("I72").Value = 13000
$ws.Range("K72").Value = 117000
$ws.Range("M72").Value = -112632

# Row 80
$ws.Range("H80").Value = 585.3333
$ws.Range("I80").Value = 371.83334
$ws.Range("J80").Value = 798.8333
$ws.Range("K80").Value = 1115.50002
$ws.Range("L80").Value = 2396.4999
$ws.Range("M80").Value = -117.5000199999999
$ws.Range("N80").Value = -4392.4999

# Row 83
$ws.Range("H83").Value = 585.3333
$ws.Range("I83").Value = 371.83334
$ws.Range("J83").Value = 798.8333
$ws.Range("K83").Value = 3346.50006
$ws.Range("L83").Value = 7189.4997
$ws.Range("M83").Value = 1645.49994
$ws.Range("N83").Value = -17173.4997

# Row 113
$ws.Range("H113").Value = 1788.6
$ws.Range("I113").Value = 1485.75
$ws.Range("J113").Value = 3000
$ws.Range("K113").Value = 1485.75
$ws.Range("L113").Value = 3000
$ws.Range("M113").Value = 1768.25
$ws.Range("N113").Value = -9508

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 6705.968
$ws.Range("I32").Value = 3884.6667
$ws.Range("J32").Value = 25749.75
$ws.Range("K32").Value = 3884.6667
$ws.Range("L32").Value = 25749.75
$ws.Range("M32").Value = -3597.6667
$ws.Range("N32").Value = -26323.75

# Row 122
$ws.Range("H122").Value = 2555.1667
$ws.Range("I122").Value = 2026.2
$ws.Range("J122").Value = 5200
$ws.Range("K122").Value = 6078.6
$ws.Range("L122").Value = 15600
$ws.Range("M122").Value = -3628.6
$ws.Range("N122").Value = -20500

$ws = $wb.Worksheets.Item("CRP")
# Row 7
$ws.Range("H7").Value = 77.90909000000001
$ws.Range("I7").Value = 60.857143
$ws.Range("J7").Value = 107.75
$ws.Range("K7").Value = 60.857143
$ws.Range("L7").Value = 107.75
$ws.Range("M7").Value = 52.142857
$ws.Range("N7").Value = -333.75

# Row 99
$ws.Range("H99").Value = 2811.75
$ws.Range("I99").Value = 2467.5
$ws.Range("J99").Value = 3500.25
$ws.Range("K99").Value = 2467.5
$ws.Range("L99").Value = 3500.25
$ws.Range("M99").Value = -969.5
$ws.Range("N99").Value = -6496.25

# Row 126
$ws.Range("H126").Value = 2811.75
$ws.Range("I126").Value = 2467.5
$ws.Range("J126").Value = 3500.25
$ws.Range("K126").Value = 7402.5
$ws.Range("L126").Value = 10500.75
$ws.Range("M126").Value = -4932.5
$ws.Range("N126").Value = -15440.75

# Row 132
$ws.Range("H132").Value = 4803.25
$ws.Range("I132").Value = 4199.2856
$ws.Range("K132").Value = 12597.8568
$ws.Range("M132").Value = -10067.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 6
$ws.Range("H6").Value = 68.77778000000001
$ws.Range("I6").Value = 79.28570999999999
$ws.Range("K6").Value = 237.85713
$ws.Range("M6").Value = -124.85713

# Row 44
$ws.Range("H44").Value = 962.3
$ws.Range("I44").Value = 270.5
$ws.Range("J44").Value = 2000
$ws.Range("K44").Value = 811.5
$ws.Range("L44").Value = 6000
$ws.Range("M44").Value = -413.5
$ws.Range("N44").Value = -6796

$ws = $wb.Worksheets.Item("GSM")
# Row 15
$ws.Range("H15").Value = 35998
$ws.Range("J15").Value = 35998
$ws.Range("L15").Value = 35998
$ws.Range("N15").Value = -36574

# Row 80
$ws.Range("H80").Value = 1684.375
$ws.Range("I80").Value = 1852.7142
$ws.Range("J80").Value = 506
$ws.Range("K80").Value = 1852.7142
$ws.Range("L80").Value = 506
$ws.Range("M80").Value = -854.7141999999999
$ws.Range("N80").Value = -2502

# Row 81
$ws.Range("H81").Value = 35998
$ws.Range("J81").Value = 35998
$ws.Range("L81").Value = 35998
$ws.Range("N81").Value = -37994

# Row 83
$ws.Range("H83").Value = 1684.375
$ws.Range("I83").Value = 1852.7142
$ws.Range("J83").Value = 506
$ws.Range("K83").Value = 9263.571
$ws.Range("L83").Value = 2530
$ws.Range("M83").Value = -4271.571
$ws.Range("N83").Value = -12514

# Row 84
$ws.Range("H84").Value = 35998
$ws.Range("J84").Value = 35998
$ws.Range("L84").Value = 107994
$ws.Range("N84").Value = -117978

# Row 102
$ws.Range("H102").Value = 1273
$ws.Range("J102").Value = 244.5
$ws.Range("L102").Value = 244.5
$ws.Range("N102").Value = -3488.5

# Row 107
$ws.Range("H107").Value = 227.45454
$ws.Range("I107").Value = 142
$ws.Range("K107").Value = 142
$ws.Range("M107").Value = 1778

# Row 122
$ws.Range("H122").Value = 7225.75
$ws.Range("J122").Value = 7225.75
$ws.Range("L122").Value = 21677.25
$ws.Range("N122").Value = -26577.25

# Row 126
$ws.Range("H126").Value = 6829.7144
$ws.Range("I126").Value = 6558.8
$ws.Range("K126").Value = 19676.4
$ws.Range("M126").Value = -17206.4

$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 6799.222
$ws.Range("I7").Value = 5074.5
$ws.Range("K7").Value = 5074.5
$ws.Range("M7").Value = -4962.5

# Row 62
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = ""

# Row 65
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = ""

# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").Value = ""

# Row 122
$ws.Range("H122").Value = 2986.4285
$ws.Range("I122").Value = 2928.6667
$ws.Range("K122").Value = 8786.000100000001
$ws.Range("M122").Value = -6336.000100000001

# Row 126
$ws.Range("H126").Value = 6799.222
$ws.Range("I126").Value = 5074.5
$ws.Range("K126").Value = 15223.5
$ws.Range("M126").Value = -12753.5

# Row 136
$ws.Range("H136").Value = 2300
$ws.Range("I136").Value = 2100
$ws.Range("J136").Value = 2500
$ws.Range("K136").Value = 6300
$ws.Range("L136").Value = 7500
$ws.Range("M136").Value = -3750
$ws.Range("N136").Value = -12600

$ws = $wb.Worksheets.Item("WVR")
# Row 7
$ws.Range("H7").Value = 1500.5
$ws.Range("J7").Value = 1000
$ws.Range("L7").Value = 1000
$ws.Range("N7").Value = -1226

# Row 132
$ws.Range("H132").Value = 3500
$ws.Range("I132").Value = 3500
$ws.Range("K132").Value = 10500
$ws.Range("M132").Value = -7970
